# Auto-generated Excel COM-interop script to apply numeric value updates
# to the Chocobo_Profits leve-profit tracking workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 935.75
$ws.Range("I15").Value = 935.75
$ws.Range("K15").Value = 2807.25
$ws.Range("M15").Value = -2638.25
$ws.Range("H17").Value = 1285.8334
$ws.Range("I17").Value = 1444.4445
$ws.Range("J17").Value = 1190.6666
$ws.Range("K17").Value = 4333.333500000001
$ws.Range("L17").Value = 3571.9998
$ws.Range("M17").Value = -4165.333500000001
$ws.Range("N17").Value = -3907.9998
$ws.Range("H18").Value = 465
$ws.Range("I18").Value = 190
$ws.Range("J18").Value = 556.6667
$ws.Range("K18").Value = 190
$ws.Range("L18").Value = 556.6667
$ws.Range("M18").Value = 94
$ws.Range("N18").Value = -1124.6667
$ws.Range("H40").Value = 1300
$ws.Range("I40").Value = 1200
$ws.Range("J40").Value = 1400
$ws.Range("K40").Value = 1200
$ws.Range("L40").Value = 1400
$ws.Range("M40").Value = -1025
$ws.Range("N40").Value = -1750
$ws.Range("H43").Value = 883.05
$ws.Range("I43").Value = 951.44446
$ws.Range("J43").Value = 827.0909
$ws.Range("K43").Value = 951.44446
$ws.Range("L43").Value = 827.0909
$ws.Range("M43").Value = -882.44446
$ws.Range("N43").Value = -965.0909
$ws.Range("H76").Value = 3120.9167
$ws.Range("I76").Value = 3104.8572
$ws.Range("J76").Value = 3233.3333
$ws.Range("K76").Value = 3104.8572
$ws.Range("L76").Value = 3233.3333
$ws.Range("M76").Value = -2789.8572
$ws.Range("N76").Value = -3863.3333
$ws.Range("H79").Value = 3120.9167
$ws.Range("I79").Value = 3104.8572
$ws.Range("J79").Value = 3233.3333
$ws.Range("K79").Value = 3104.8572
$ws.Range("L79").Value = 3233.3333
$ws.Range("M79").Value = -2012.8572
$ws.Range("N79").Value = -5417.3333
$ws.Range("H131").Value = 3055.238
$ws.Range("I131").Value = 1431.25
$ws.Range("J131").Value = 5220.5557
$ws.Range("K131").Value = 4293.75
$ws.Range("L131").Value = 15661.6671
$ws.Range("M131").Value = 746.25
$ws.Range("N131").Value = -25741.6671
$ws.Range("H132").Value = 5395.1284
$ws.Range("I132").Value = 5137.4688
$ws.Range("J132").Value = 6573
$ws.Range("K132").Value = 15412.4064
$ws.Range("L132").Value = 19719
$ws.Range("M132").Value = -12882.4064
$ws.Range("N132").Value = -24779
$ws.Range("H137").Value = 3443.158
$ws.Range("I137").Value = 1778
$ws.Range("J137").Value = 9687.5
$ws.Range("K137").Value = 5334
$ws.Range("L137").Value = 29062.5
$ws.Range("M137").Value = -2784
$ws.Range("N137").Value = -34162.5
$ws.Range("H138").Value = 3300.658
$ws.Range("I138").Value = 1713.3846
$ws.Range("J138").Value = 3628.1904
$ws.Range("K138").Value = 5140.1538
$ws.Range("L138").Value = 10884.5712
$ws.Range("M138").Value = -0.1538000000000466
$ws.Range("N138").Value = -21164.5712

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 26542.777
$ws.Range("I37").Value = 29800
$ws.Range("J37").Value = 26135.625
$ws.Range("K37").Value = 29800
$ws.Range("L37").Value = 26135.625
$ws.Range("M37").Value = -29527
$ws.Range("N37").Value = -26681.625
$ws.Range("H74").Value = 4175.7666
$ws.Range("I74").Value = 4704.048
$ws.Range("J74").Value = 2943.111
$ws.Range("K74").Value = 4704.048
$ws.Range("L74").Value = 2943.111
$ws.Range("M74").Value = -3830.048
$ws.Range("N74").Value = -4691.111
$ws.Range("H77").Value = 4175.7666
$ws.Range("I77").Value = 4704.048
$ws.Range("J77").Value = 2943.111
$ws.Range("K77").Value = 23520.24
$ws.Range("L77").Value = 14715.555
$ws.Range("M77").Value = -19152.24
$ws.Range("N77").Value = -23451.555

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 30807.555
$ws.Range("J35").Value = 30807.555
$ws.Range("L35").Value = 30807.555
$ws.Range("N35").Value = -31427.555
$ws.Range("H124").Value = 41780
$ws.Range("J124").Value = 41780
$ws.Range("L124").Value = 41780
$ws.Range("N124").Value = -51600
$ws.Range("H140").Value = 51319
$ws.Range("J140").Value = 51319
$ws.Range("L140").Value = 51319
$ws.Range("N140").Value = -61679

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7938535.5
$ws.Range("I16").Value = 13890060
$ws.Range("J16").Value = 3168.8333
$ws.Range("K16").Value = 13890060
$ws.Range("L16").Value = 3168.8333
$ws.Range("M16").Value = -13889773
$ws.Range("N16").Value = -3742.8333
$ws.Range("H41").Value = 30870.908
$ws.Range("J41").Value = 30870.908
$ws.Range("L41").Value = 30870.908
$ws.Range("N41").Value = -31726.908
$ws.Range("H86").Value = 2005
$ws.Range("I86").Value = 2003.5
$ws.Range("J86").Value = 2008
$ws.Range("K86").Value = 2003.5
$ws.Range("L86").Value = 2008
$ws.Range("M86").Value = -880.5
$ws.Range("N86").Value = -4254
$ws.Range("H89").Value = 2005
$ws.Range("I89").Value = 2003.5
$ws.Range("J89").Value = 2008
$ws.Range("K89").Value = 10017.5
$ws.Range("L89").Value = 10040
$ws.Range("M89").Value = -4401.5
$ws.Range("N89").Value = -21272
$ws.Range("H113").Value = 7938535.5
$ws.Range("I113").Value = 13890060
$ws.Range("J113").Value = 3168.8333
$ws.Range("K113").Value = 13890060
$ws.Range("L113").Value = 3168.8333
$ws.Range("M113").Value = -13887890
$ws.Range("N113").Value = -7508.8333
$ws.Range("H132").Value = 3186.5625
$ws.Range("I132").Value = 1744.2727
$ws.Range("J132").Value = 6359.6
$ws.Range("K132").Value = 5232.8181
$ws.Range("L132").Value = 19078.8
$ws.Range("M132").Value = -2702.8181
$ws.Range("N132").Value = -24138.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 2137
$ws.Range("I21").Value = 546.5
$ws.Range("K21").Value = 1639.5
$ws.Range("M21").Value = -1466.5
$ws.Range("H34").Value = 8177.6924
$ws.Range("I34").Value = 170.5
$ws.Range("J34").Value = 11736.444
$ws.Range("K34").Value = 511.5
$ws.Range("L34").Value = 35209.33199999999
$ws.Range("M34").Value = -427.5
$ws.Range("N34").Value = -35377.33199999999
$ws.Range("H39").Value = 14000
$ws.Range("I39").Value = 8000
$ws.Range("J39").Value = 15000
$ws.Range("K39").Value = 24000
$ws.Range("L39").Value = 45000
$ws.Range("M39").Value = -23706
$ws.Range("N39").Value = -45588
$ws.Range("H55").Value = 745
$ws.Range("I55").Value = 745
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 2235
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -2058
$ws.Range("N55").ClearContents()
$ws.Range("H68").Value = 9953.182000000001
$ws.Range("I68").Value = 616
$ws.Range("J68").Value = 17734.166
$ws.Range("K68").Value = 1848
$ws.Range("L68").Value = 53202.49800000001
$ws.Range("M68").Value = -1037
$ws.Range("N68").Value = -54824.49800000001
$ws.Range("H71").Value = 9953.182000000001
$ws.Range("I71").Value = 616
$ws.Range("J71").Value = 17734.166
$ws.Range("K71").Value = 5544
$ws.Range("L71").Value = 159607.494
$ws.Range("M71").Value = -1488
$ws.Range("N71").Value = -167719.494
$ws.Range("H113").Value = 793.8378
$ws.Range("I113").Value = 662.5217
$ws.Range("K113").Value = 1987.5651
$ws.Range("M113").Value = 182.4349
$ws.Range("H139").Value = 3031.0527
$ws.Range("I139").Value = 1512.8572
$ws.Range("J139").Value = 3916.6667
$ws.Range("K139").Value = 4538.571599999999
$ws.Range("L139").Value = 11750.0001
$ws.Range("M139").Value = 601.4284000000007
$ws.Range("N139").Value = -22030.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2474.889
$ws.Range("I102").Value = 1660.875
$ws.Range("K102").Value = 1660.875
$ws.Range("M102").Value = -38.875
$ws.Range("H132").Value = 6085.2856
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 25500
$ws.Range("J18").Value = 25500
$ws.Range("L18").Value = 25500
$ws.Range("M18").Value = -25844
$ws.Range("H127").Value = 30174.445
$ws.Range("J127").Value = 30174.445
$ws.Range("L127").Value = 30174.445
$ws.Range("N127").Value = -40094.445
